$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diagram is gaining one more level of nesting for the Farzi -> Bike
# dependency arrow, so every existing column from B to G needs to shift one
# column to the right (B->C, C->D, D->E, E->F, F->G, G->H), leaving column A
# (the outermost "Bike"/"Car" branch indicator) untouched.
$ws.Columns("B").Insert()

# Inserting a whole column stamps every row in the new column B with the
# bordered style inherited from its neighbours. That vertical connector line
# is only actually part of the diagram for rows 18-21 (the new
# Farzi -> Bike dependency block), so clear it everywhere else to match the
# original hand-built layout.
$ws.Range("B2:B17").Clear()

# The dependency/inheritance arrows are being re-drawn with a heavier red
# double-arrow glyph instead of the thin single-arrow glyph used before.
[void]$ws.UsedRange.Replace("$([char]0x2192)", "$([char]0x21D2)")

# Draw the new arrow/triangle glyphs for the extra Farzi -> Bike dependency
# hop that is now needed after the re-indent above.
$ws.Range("B18").Value = "$([char]0x21D2)"
$ws.Range("G18").Value = "$([char]0x25B7)"
$ws.Range("B21").Value = "$([char]0x25C1)"
